$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.0001195504850560666
$ws.Range("E2").Value = 0.02111284914371071
$ws.Range("D3").Value = 0.1745105794506428
$ws.Range("E3").Value = 1.465313996513155
$ws.Range("D4").Value = 0.03742397646105173
$ws.Range("E4").Value = 0.369150964590698
$ws.Range("D5").Value = -0.8269341827440561
$ws.Range("E5").Value = 2.719615507367969
$ws.Range("D6").Value = -0.002198448813783861
$ws.Range("E6").Value = 0.01739034962636656
$ws.Range("D7").Value = 0.852517565939068
$ws.Range("E7").Value = 4.513242897748221
$ws.Range("D8").Value = 0.09357697738521245
$ws.Range("E8").Value = 0.0880883822478256
$ws.Range("D9").Value = 0.3053119412037891
$ws.Range("E9").Value = 4.656063420286944
$ws.Range("D10").Value = -0.0005182812517577343
$ws.Range("E10").Value = 0.04772616952074563
$ws.Range("D11").Value = 0.1983290592648199
$ws.Range("E11").Value = 3.584565954788263
$ws.Range("D12").Value = -0.1256624957323176
$ws.Range("E12").Value = 0.4139875257300115
$ws.Range("D13").Value = -0.1244015628988562
$ws.Range("E13").Value = 4.399578731105592
$ws.Range("D14").Value = 0.0007842094692209424
$ws.Range("E14").Value = 0.01915594907124332
$ws.Range("D15").Value = 0.04748792940276472
$ws.Range("E15").Value = 1.132811215694494
$ws.Range("D16").Value = 0.1680308445848331
$ws.Range("E16").Value = 0.05630954063168365
$ws.Range("D17").Value = 2.51856664492156
$ws.Range("E17").Value = 1.882411488095351
